$wb = $excel.ActiveWorkbook

# --- Sheet "展览" (sheetId=1) ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Cells.Item(2, 6).Value  = 6526   # F2  6516 -> 6526
$ws1.Cells.Item(3, 6).Value  = 118    # F3  115  -> 118
$ws1.Cells.Item(5, 6).Value  = 403    # F5  401  -> 403
$ws1.Cells.Item(8, 6).Value  = 537    # F8  536  -> 537
$ws1.Cells.Item(9, 6).Value  = 93     # F9  92   -> 93
$ws1.Cells.Item(13, 6).Value = 382    # F13 383  -> 382
$ws1.Cells.Item(14, 6).Value = 957    # F14 953  -> 957
$ws1.Cells.Item(15, 6).Value = 3225   # F15 3221 -> 3225
$ws1.Cells.Item(16, 7).Value = 30     # G16 45   -> 30
$ws1.Cells.Item(18, 6).Value = 1881   # F18 1880 -> 1881

# --- Sheet "全部类型" (sheetId=4) ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Cells.Item(2, 6).Value  = 6526   # F2  6516 -> 6526
$ws4.Cells.Item(3, 6).Value  = 118    # F3  115  -> 118
$ws4.Cells.Item(5, 6).Value  = 403    # F5  401  -> 403
$ws4.Cells.Item(9, 6).Value  = 537    # F9  536  -> 537
$ws4.Cells.Item(10, 6).Value = 93     # F10 92   -> 93
$ws4.Cells.Item(15, 6).Value = 957    # F15 953  -> 957
$ws4.Cells.Item(16, 6).Value = 3225   # F16 3221 -> 3225
$ws4.Cells.Item(17, 7).Value = 30     # G17 45   -> 30
$ws4.Cells.Item(19, 6).Value = 1881   # F19 1880 -> 1881
